# Insert a new weekly price-report row for Perejil (Vega Modelo de Temuco) at
# row 293, pushing every following row (old 293..424) down by one, which also
# makes the sheet's used range grow from A1:R424 to A1:R425 automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(293).Insert()

$ws.Range("A293").Value = 10
$ws.Range("B293").Value = "Vega Modelo de Temuco"
$ws.Range("C293").Value = "La Araucanía"
$ws.Range("D293").Value = 44917
$ws.Range("E293").Value = 9
$ws.Range("F293").Value = 100112044
$ws.Range("G293").Value = "Perejil"
$ws.Range("H293").Value = "Sin especificar"
$ws.Range("I293").Value = "Primera"
$ws.Range("J293").Value = 65
$ws.Range("K293").Value = 5000
$ws.Range("L293").Value = 5000
$ws.Range("M293").Value = 5000
$ws.Range("N293").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O293").Value = "Provincia de Cautín"
$ws.Range("P293").Value = 1667
$ws.Range("Q293").Value = 3
$ws.Range("R293").Value = "Hortaliza"
